$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.949.04"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.67%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.985.52"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.19%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.06%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.73"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -4.75%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.06"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("E9").Value = "  -4.45%  "
$ws.Range("E10").Value = "  -4.33%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.359"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.74%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.496.03"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("E13").Value = "  -2.28%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.12"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("E15").Value = "  -7.17%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.053.30"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.42%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.08"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -2.09%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.55"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -2.08%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E20").Value = "  -3.89%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.64"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -5.85%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = $style
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.76%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.492"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -1.59%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.64"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("E26").Value = "  +0.17%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -5.09%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0892"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -8.10%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.55"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -6.62%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("E32").Value = "  -6.78%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.23"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -4.16%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.76"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("E36").Value = "  -1.83%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -6.66%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.31"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -6.04%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0666"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -4.49%  "
$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.016.15"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.70"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +0.00%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -3.01%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.643"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.92%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.204.25"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -5.37%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.39"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -5.87%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.95"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.942"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -8.93%  "
$ws.Range("E49").Value = "  -5.36%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.31"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("E51").Value = "  -11.61%  "
